$wb = $excel.ActiveWorkbook

# 1. Insert the new "manufacturer_country" worksheet right after "manufacturer"
$mfgSheet = $wb.Worksheets.Item("manufacturer")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $mfgSheet)
$ws.Name = "manufacturer_country"

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Original name"
$ws.Range("C1").Value = "Data Type"
$ws.Range("D1").Value = "Length"
$ws.Range("E1").Value = "Nullable"
$ws.Range("F1").Value = "Default"
$ws.Range("G1").Value = "Description"
$ws.Range("A1:G1").Font.Bold = $true

$ws.Range("A2").Value = "manufacturerID"
$ws.Range("C2").Value = "smallint"
$ws.Range("E2").Value = "NO"

$ws.Range("A3").Value = "countryID"
$ws.Range("C3").Value = "int"
$ws.Range("E3").Value = "NO"

$ws.Range("A5").Value = "References (Child)"
$ws.Range("A5").Font.Bold = $true
$ws.Range("B5").Value = "This table references the following parent tables:"

$ws.Range("A6").Value = "manufacturer"
$ws.Range("A7").Value = "country"

Write-Output "done"
